$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the TestResult (column G) values for the data rows (2-23), keeping formatting.
$ws.Range("G2:G23").ClearContents()

# Update the active selection to match the new cursor position.
$ws.Range("F17").Select()
